$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Euramet")

# Swap F7 and F9 values
$ws.Range("F7").Value = -156.5579745837849
$ws.Range("F9").Value = -157.0411782090434

# Row 12: only E changes
$ws.Range("E12").Value = 0

# Rows 13-15: E -> 0, F swapped to the other value
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = -156.5579745837849

$ws.Range("E14").Value = 0
$ws.Range("F14").Value = -156.5579745837849

$ws.Range("E15").Value = 0
$ws.Range("F15").Value = -156.5579745837849

# Rows 16-19: clear D:H contents entirely (keep formatting)
$ws.Range("D16:H19").ClearContents()

# Swap F29 and F32/F33
$ws.Range("F29").Value = -156.5579745837849
$ws.Range("F32").Value = -157.0411782090434
$ws.Range("F33").Value = -157.0411782090434

# Rows 38-41: clear D:H contents entirely (keep formatting)
$ws.Range("D38:H41").ClearContents()
